$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 172, shifting existing rows 172:194 down to 173:195.
# Excel's default Insert() copies formatting from the row above, which matches
# the target column D style (s="2") on the new row.
$ws.Rows.Item(172).Insert()

# Populate the newly inserted row 172 with the new record's data.
$ws.Range("A172").Value = 5
$ws.Range("B172").Value = "Macroferia Regional de Talca"
$ws.Range("C172").Value = "Maule"
$ws.Range("D172").Value = 45142
$ws.Range("E172").Value = 7
$ws.Range("F172").Value = 100112001
$ws.Range("G172").Value = "Berenjena"
$ws.Range("H172").Value = "Sin especificar"
$ws.Range("I172").Value = "Primera"
$ws.Range("J172").Value = 150
$ws.Range("K172").Value = 7000
$ws.Range("L172").Value = 7000
$ws.Range("M172").Value = 7000
$ws.Range("N172").Value = "$/caja 50 unidades"
$ws.Range("O172").Value = "Región de Arica y Parinacota"
$ws.Range("P172").Value = 140
$ws.Range("Q172").Value = 50
$ws.Range("R172").Value = "Hortaliza"
